# Updates "想去人数" (interest count) values in the 展览 (Exhibition),
# 演出 (Performance) and 全部类型 (All Types) sheets, matching the
# gh-pages data regeneration at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1016
$ws1.Range("F4").Value = 10
$ws1.Range("F5").Value = 7497
$ws1.Range("F7").Value = 927
$ws1.Range("F8").Value = 285
$ws1.Range("F11").Value = 81
$ws1.Range("F15").Value = 3006
$ws1.Range("F16").Value = 175
$ws1.Range("F17").Value = 72
$ws1.Range("F18").Value = 699
$ws1.Range("F19").Value = 746
$ws1.Range("F21").Value = 439
$ws1.Range("F23").Value = 186
$ws1.Range("F24").Value = 205
$ws1.Range("F25").Value = 215
$ws1.Range("F26").Value = 239
$ws1.Range("F28").Value = 87
$ws1.Range("F29").Value = 225
$ws1.Range("F32").Value = 372
$ws1.Range("F33").Value = 428
$ws1.Range("F37").Value = 70

# 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 34

# 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1016
$ws4.Range("F5").Value = 34
$ws4.Range("F7").Value = 10
$ws4.Range("F8").Value = 7497
$ws4.Range("F10").Value = 927
$ws4.Range("F11").Value = 285
$ws4.Range("F14").Value = 81
$ws4.Range("F19").Value = 3006
$ws4.Range("F20").Value = 175
$ws4.Range("F21").Value = 72
$ws4.Range("F23").Value = 699
$ws4.Range("F24").Value = 746
$ws4.Range("F27").Value = 439
$ws4.Range("F29").Value = 186
$ws4.Range("F30").Value = 205
$ws4.Range("F31").Value = 215
$ws4.Range("F32").Value = 239
$ws4.Range("F34").Value = 87
$ws4.Range("F35").Value = 225
$ws4.Range("F38").Value = 372
$ws4.Range("F39").Value = 428
$ws4.Range("F43").Value = 70

